$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$np = $s.NotesPage
$np.Shapes.Placeholders.Item(2).TextFrame.TextRange.Text = "Hello notes"
Write-Host ("HasNotesMaster=" + $p.HasNotesMaster)
$nm = $p.NotesMaster
Write-Host ("HasNotesMaster2=" + $p.HasNotesMaster)
